{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the target paragraph: the one right after the \"You use the execution()...\"\n// paragraph, containing the placeholder text \"Now..\" (with a first-line indent).\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Now..\") {\n    targetIndex = i;\n    break;\n  }\n}\nif (targetIndex === -1) {\n  throw new Error(\"Could not locate the 'Now..' placeholder paragraph.\");\n}\n\n// Build the OOXML package fragment holding all of the replacement paragraphs\n// (the long passage about within(), bean() pointcut designators and annotated\n// aspects), ending with the paragraph that keeps the existing _GoBack bookmark.\nconst bodyFragment = `<w:p><w:r><w:tab/><w:t>Now let\u2019s suppose that you want  to confine</w:t></w:r><w:r><w:t xml:space=\"preserve\"> the reach of that pointcut to only the concert package. In that case, you can limit the match by tacking on a within() designator, as a shown at bottom:</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">//------: execution(* concert.Performance.perform(..)) &amp;&amp; within(concert.*) </w:t></w:r></w:p><w:p><w:r><w:t>Note that you use the &amp;&amp; operator to combine the execution() and within designators in an \u201cand\u201d relationship (where both designators must match for the pointcut to match).</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Similary, you could use the || operator to indicate an \u201cor\u201d relationship.  And the ! operator can be used to negate the effect of  a designator.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:tab/><w:t>Because ampersands have special meaning in XML-based configuration. Likewise, or and not can be used in place of || and ! , respectively.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:t>Selecting beans in pointcuts</w:t></w:r></w:p><w:p><w:r><w:t>Spring adds a bean() designator  that  lets you identify beans by their ID in a pointcut expression. bean()  takes a bean ID or names as an argument and limits the pointcut\u2019s effect to that specific bean.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>For example, consider the following pointcut:</w:t></w:r></w:p><w:p><w:r><w:t>//------execution(* concert.Performance.perform())  and bean(\u2018woodstock\u2019)</w:t></w:r></w:p><w:p><w:r><w:t>Here you\u2019re saying  that  you want  to apply aspect  advice to the execution of  Performance\u2019s perform() method, but limited to the bean whose ID is woodstock.  Narrowing  a pointcut to a specific bean may be valuable in some  cases, but you can also use negation to apply an aspect to all beans that don\u2019t have  a specific ID:</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">//------ execution(* concert.Performance.perform()) and !bean(\u2018woodstock\u2019) </w:t></w:r></w:p><w:p><w:r><w:t>In this case, the aspect\u2019s advice will be woven into all beans whose ID isn\u2019t  \u2018woodstock\u2019.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:t>Creating annotated aspects</w:t></w:r></w:p><w:p><w:r><w:t>A key feature introduced  in  AspectJ 5 is the ability to use annotati</w:t></w:r><w:r><w:t xml:space=\"preserve\">ons  to create aspects. Prior to </w:t></w:r><w:r><w:t>AspectJ 5, wiring AspectJ aspects involved learning a Java language extension. But AspectJ 5, wiring AspectJ\u2019s annotation-oriented model makes it simple to turn any class into an aspect by sprinkling a few annotations around.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>You\u2019ve already defined the Performance interface as the subject of  your  aspect\u2019s pointcuts.Now let\u2019s use AspectJ annotations to create an aspect.</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>`;\nconst ooxmlPackage = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + bodyFragment + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nconst targetParagraph = paragraphs.items[targetIndex];\nconst targetRange = targetParagraph.getRange();\n// Replacing this single paragraph's range with the multi-paragraph OOXML\n// fragment above expands it into the full set of new paragraphs, while the\n// bookmarkStart/bookmarkEnd that trailed the old paragraph are left alone\n// (they stay attached right after the inserted content).\ntargetRange.insertOoxml(ooxmlPackage, Word.InsertLocation.replace);\nawait context.sync();\n\n// The original document ended with four empty trailing paragraphs; the\n// edit trims that down to two (the very last paragraph in a Word document\n// can never be removed, so we delete the two empty ones just before it).\nconst trailingParagraphs = context.document.body.paragraphs;\ntrailingParagraphs.load(\"items,text,isLastParagraph\");\nawait context.sync();\n\nconst items = trailingParagraphs.items;\nlet emptyTrailing = [];\nfor (let i = items.length - 1; i >= 0; i--) {\n  if (items[i].text === \"\") {\n    emptyTrailing.push(items[i]);\n  } else {\n    break;\n  }\n}\n// Keep the very last paragraph and one more; delete the rest of the\n// contiguous trailing empty paragraphs so exactly two remain.\nconst toDelete = emptyTrailing.slice(2);\nfor (const para of toDelete) {\n  para.delete();\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the target paragraph: the one right after the \"You use the execution()...\"\n# paragraph, containing the placeholder text \"Now..\" (with a first-line indent).\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq \"Now..\") {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the 'Now..' placeholder paragraph.\"\n}\n\n# Build the OOXML package fragment holding all of the replacement paragraphs\n# (the long passage about within(), bean() pointcut designators and annotated\n# aspects), ending with the paragraph that keeps the existing _GoBack bookmark.\n$bodyFragment = '<w:p><w:r><w:tab/><w:t>Now let\u2019s suppose that you want  to confine</w:t></w:r><w:r><w:t xml:space=\"preserve\"> the reach of that pointcut to only the concert package. In that case, you can limit the match by tacking on a within() designator, as a shown at bottom:</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">//------: execution(* concert.Performance.perform(..)) &amp;&amp; within(concert.*) </w:t></w:r></w:p><w:p><w:r><w:t>Note that you use the &amp;&amp; operator to combine the execution() and within designators in an \u201cand\u201d relationship (where both designators must match for the pointcut to match).</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Similary, you could use the || operator to indicate an \u201cor\u201d relationship.  And the ! operator can be used to negate the effect of  a designator.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:tab/><w:t>Because ampersands have special meaning in XML-based configuration. Likewise, or and not can be used in place of || and ! , respectively.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:t>Selecting beans in pointcuts</w:t></w:r></w:p><w:p><w:r><w:t>Spring adds a bean() designator  that  lets you identify beans by their ID in a pointcut expression. bean()  takes a bean ID or names as an argument and limits the pointcut\u2019s effect to that specific bean.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>For example, consider the following pointcut:</w:t></w:r></w:p><w:p><w:r><w:t>//------execution(* concert.Performance.perform())  and bean(\u2018woodstock\u2019)</w:t></w:r></w:p><w:p><w:r><w:t>Here you\u2019re saying  that  you want  to apply aspect  advice to the execution of  Performance\u2019s perform() method, but limited to the bean whose ID is woodstock.  Narrowing  a pointcut to a specific bean may be valuable in some  cases, but you can also use negation to apply an aspect to all beans that don\u2019t have  a specific ID:</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">//------ execution(* concert.Performance.perform()) and !bean(\u2018woodstock\u2019) </w:t></w:r></w:p><w:p><w:r><w:t>In this case, the aspect\u2019s advice will be woven into all beans whose ID isn\u2019t  \u2018woodstock\u2019.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:t>Creating annotated aspects</w:t></w:r></w:p><w:p><w:r><w:t>A key feature introduced  in  AspectJ 5 is the ability to use annotati</w:t></w:r><w:r><w:t xml:space=\"preserve\">ons  to create aspects. Prior to </w:t></w:r><w:r><w:t>AspectJ 5, wiring AspectJ aspects involved learning a Java language extension. But AspectJ 5, wiring AspectJ\u2019s annotation-oriented model makes it simple to turn any class into an aspect by sprinkling a few annotations around.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>You\u2019ve already defined the Performance interface as the subject of  your  aspect\u2019s pointcuts.Now let\u2019s use AspectJ annotations to create an aspect.</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>'\n$ooxmlPackage = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + $bodyFragment + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n# Replacing this single paragraph's Range with the multi-paragraph OOXML\n# fragment above expands it into the full set of new paragraphs. The\n# paragraph's Range includes the trailing bookmarkStart/bookmarkEnd for\n# _GoBack, so the fragment re-creates that bookmark on the final paragraph.\n$targetRange = $d.Paragraphs.Item($targetIndex).Range\n$targetRange.InsertXML($ooxmlPackage)\n\n# The original document ended with four empty trailing paragraphs; the\n# edit trims that down to two (the very last paragraph in a Word document\n# can never be removed, so we delete the two empty ones just before it).\n$emptyTrailing = @()\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq \"\") {\n        $emptyTrailing += $i\n    } else {\n        break\n    }\n}\n# $emptyTrailing is ordered from the last paragraph backwards; keep the\n# last two (closest to the end) and delete the rest.\nif ($emptyTrailing.Count -gt 2) {\n    $toDelete = $emptyTrailing[2..($emptyTrailing.Count - 1)]\n    foreach ($idx in $toDelete) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
